$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column D (which currently holds "Tipo").
# This shifts the existing "Tipo"/"multiple" column from D to E.
$ws.Columns("D").Insert()

# Set the new header "MAE" in D1 and copy the header formatting
# (bold font, border, centered alignment) from the other header cells.
$ws.Range("D1").Value = "MAE"
$ws.Range("A1").Copy()
$ws.Range("D1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Fill in the MAE values for each data row
$ws.Range("D2").Value = 1.027519701295929
$ws.Range("D3").Value = 1.350102492671558
$ws.Range("D4").Value = 2.046624516121589
$ws.Range("D5").Value = 2.091210877959376
